$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price updates in column D (stored as text to match source data format).
# Each cell is kept text-formatted ("@") before assignment so that
# numeric-looking values (e.g. "242.77") are not coerced to Double by COM.
$priceUpdates = @{
    "D2" = "242.77"
    "D3" = "23.10"
    "D4" = "5.411"
    "D7" = "6.524"
    "D8" = "0.8089"
    "D9" = "0.9290"
    "D10" = "0.1419"
    "D11" = "0.07387"
    "D12" = "0.03293"
    "D13" = "0.03060"
    "D14" = "0.09366"
    "D15" = "3.858"
    "D16" = "0.001585"
    "D17" = "0.04668"
    "D18" = "0.0005941"
    "D19" = "0.005928"
    "D20" = "0.001265"
    "D21" = "0.004904"
    "D24" = "2.147"
    "D26" = "0.1331"
    "D40" = "0.03974"
    "D41" = "0.006176"
    "D43" = "0.002570"
    "D44" = "0.008757"
    "D45" = "0.00005183"
    "D48" = "0.002388"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

Write-Host "Updated $($priceUpdates.Count) price cells"